$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell values (text rearrangement) ---
$ws.Cells.Item(8,1).Value = "Weather effects e.g. rain"
$ws.Cells.Item(9,1).Value = "Enemy spawning animation"
$ws.Cells.Item(9,2).Value = "Wave-intermission music"
$ws.Cells.Item(10,1).Value = "Walking bobbing"
$ws.Cells.Item(12,2).Value = "Have lumbering enemy giants that are slow much take lots of hits"

# --- Add the new "implementation status" legend rows ---
$ws.Cells.Item(24,1).Value = "Fully implemented"
$ws.Cells.Item(25,1).Value = "Partially implemented"
$ws.Cells.Item(26,1).Value = "Not implemented"

# --- New feature ideas added to the list ---
$ws.Cells.Item(10,2).Value = "Sprinting with Shift"
$ws.Cells.Item(11,2).Value = "Create a sniper"

# --- Apply "Good" / "Neutral" cell styles to reflect implementation status ---
$ws.Cells.Item(3,2).Style = "Good"
$ws.Cells.Item(3,3).Style = "Neutral"
$ws.Cells.Item(4,2).Style = "Neutral"
$ws.Cells.Item(6,2).Style = "Good"
$ws.Cells.Item(7,3).Style = "Neutral"
$ws.Cells.Item(9,3).Style = "Neutral"
$ws.Cells.Item(10,2).Style = "Good"
$ws.Cells.Item(11,3).Style = "Good"
$ws.Cells.Item(24,1).Style = "Good"
$ws.Cells.Item(25,1).Style = "Neutral"

# --- Update selection to reflect where the author ended up ---
$ws.Range("A12").Select() | Out-Null
